$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.57"
$ws.Range("D3").Value = "'22.99"
$ws.Range("D4").Value = "'5.416"
$ws.Range("D5").Value = "'0.05941"
$ws.Range("D7").Value = "'6.546"
$ws.Range("D8").Value = "'0.8125"
$ws.Range("D9").Value = "'0.9168"
$ws.Range("D11").Value = "'0.07489"
$ws.Range("D12").Value = "'0.03292"
$ws.Range("D14").Value = "'0.09351"
$ws.Range("D15").Value = "'3.849"
$ws.Range("D16").Value = "'0.001556"
$ws.Range("D18").Value = "'0.0005938"
$ws.Range("D19").Value = "'0.006081"
$ws.Range("D21").Value = "'0.0009857"
$ws.Range("E22").Value = "21NitroExNTXBestin24h"
$ws.Range("D24").Value = "'2.150"
$ws.Range("D26").Value = "'0.1323"
$ws.Range("D40").Value = "'0.03955"
$ws.Range("D41").Value = "'0.006218"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1074"
$ws.Range("D43").Value = "'0.002620"
$ws.Range("D44").Value = "'0.008935"
$ws.Range("D45").Value = "'0.00005228"
$ws.Range("D49").Value = "'0.002268"
